$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume snapshot (GitHub Actions data pull).
# Coin names, links, and Volume(1h) percentages are plain text - assign directly.
$ws.Range("D2").Value = '52.093.57'
$ws.Range("E2").Value = '  +5.40%  '
$ws.Range("D3").Value = '2.790.19'
$ws.Range("E3").Value = '  +6.33%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +5.34%  '
$ws.Range("E6").Value = '  +5.66%  '
$ws.Range("E7").Value = '  +3.79%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +6.98%  '
$ws.Range("E10").Value = '  +7.92%  '
$ws.Range("E11").Value = '  +7.61%  '
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("E13").Value = '  +2.50%  '
$ws.Range("E14").Value = '  +4.50%  '
$ws.Range("D15").Value = '3.232.63'
$ws.Range("E15").Value = '  +6.44%  '
$ws.Range("D16").Value = '2.805.56'
$ws.Range("E16").Value = '  +6.52%  '
$ws.Range("E17").Value = '  +4.67%  '
$ws.Range("D18").Value = '51.979.43'
$ws.Range("E18").Value = '  +5.30%  '
$ws.Range("E19").Value = '  +12.26%  '
$ws.Range("E20").Value = '  +3.65%  '
$ws.Range("E21").Value = '  +4.97%  '
$ws.Range("D22").Value = '0.0₃0984'
$ws.Range("E22").Value = '  +4.28%  '
$ws.Range("E23").Value = '  +4.26%  '
$ws.Range("E24").Value = '  +2.37%  '
$ws.Range("E25").Value = '  +10.58%  '
$ws.Range("E26").Value = '  +3.71%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("E30").Value = '  +3.91%  '
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("E32").Value = '  +1.51%  '
$ws.Range("E33").Value = '  +4.67%  '
$ws.Range("E34").Value = '  +2.21%  '
$ws.Range("E35").Value = '  +6.04%  '
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("E38").Value = '  +7.60%  '
$ws.Range("E39").Value = '  +1.43%  '
$ws.Range("E40").Value = '  +26.74%  '
$ws.Range("E41").Value = '  +14.37%  '
$ws.Range("E42").Value = '  +4.59%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E43").Value = '  -0.97%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E44").Value = '  +4.71%  '
$ws.Range("E45").Value = '  +3.81%  '
$ws.Range("D46").Value = '2.108.33'
$ws.Range("E46").Value = '  +3.52%  '
$ws.Range("E47").Value = '  +4.72%  '
$ws.Range("E48").Value = '  +3.38%  '
$ws.Range("E49").Value = '  +7.35%  '
$ws.Range("E50").Value = '  +21.52%  '
$ws.Range("E51").Value = '  +1.32%  '

# Price cells whose new text happens to parse as a plain number need special
# handling: a direct .Value assignment would silently convert the string to a
# real number (dropping formatting such as trailing zeros, e.g. "20.10" -> 20.1,
# "1.00" -> 1), which the source data does not do (Price is stored as text).
# Writing it as a quoted-text formula and then Paste-Special-ing (values only)
# over itself keeps the exact literal text without leaving a formula behind and
# without altering the cells formatting/style.
$ws.Range("D5").Formula = '="116.78"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = '="343.19"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D9").Formula = '="0.580"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("D10").Formula = '="42.37"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("D11").Formula = '="0.0868"'
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("D12").Formula = '="20.10"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("D17").Formula = '="0.886"'
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("D20").Formula = '="13.39"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("D21").Formula = '="6.98"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("D23").Formula = '="277.79"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("D25").Formula = '="2.80"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("D26").Formula = '="26.87"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("D28").Formula = '="10.21"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("D30").Formula = '="0.143"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("D31").Formula = '="34.98"'
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("D32").Formula = '="50.35"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("D33").Formula = '="5.71"'
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("D36").Formula = '="1.00"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D37").Formula = '="19.02"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("D38").Formula = '="3.31"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("D39").Formula = '="4.98"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("D40").Formula = '="2.74"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("D41").Formula = '="0.0371"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("D42").Formula = '="23.70"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("D43").Formula = '="127.64"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("D44").Formula = '="2.34"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("D45").Formula = '="0.115"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("D50").Formula = '="0.911"'
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("D51").Formula = '="8.94"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$excel.CutCopyMode = $false
